$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.313.48'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').Value = '3.245.26'
$ws.Range('E3').Value = '  +3.19%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = "'595.33"
$ws.Range('E5').Value = '  -1.12%  '
$ws.Range('D6').Value = "'140.46"
$ws.Range('E6').Value = '  -1.30%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '3.241.95'
$ws.Range('E8').Value = '  +3.14%  '
$ws.Range('D9').Value = "'0.519"
$ws.Range('E9').Value = '  -1.77%  '
$ws.Range('E10').Value = '  -1.04%  '
$ws.Range('D11').Value = "'5.36"
$ws.Range('E11').Value = '  -0.61%  '
$ws.Range('D12').Value = "'0.465"
$ws.Range('E12').Value = '  -0.37%  '
$ws.Range('D13').Value = "'0.0000247"
$ws.Range('E13').Value = '  -2.70%  '
$ws.Range('D14').Value = "'34.34"
$ws.Range('E14').Value = '  -1.57%  '
$ws.Range('D15').Value = '3.778.11'
$ws.Range('E15').Value = '  +3.20%  '
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').Value = '3.240.09'
$ws.Range('E17').Value = '  +3.11%  '
$ws.Range('D18').Value = '63.352.49'
$ws.Range('E18').Value = '  -0.97%  '
$ws.Range('E19').Value = '  -1.08%  '
$ws.Range('D20').Value = "'473.99"
$ws.Range('E20').Value = '  -2.64%  '
$ws.Range('D21').Value = "'14.17"
$ws.Range('D22').Value = "'0.732"
$ws.Range('E22').Value = '  +2.89%  '
$ws.Range('D23').Value = "'7.94"
$ws.Range('E23').Value = '  +2.49%  '
$ws.Range('D24').Value = "'83.93"
$ws.Range('E25').Value = '  -0.36%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D28').Value = "'7.35"
$ws.Range('E28').Value = '  +5.38%  '
$ws.Range('E29').Value = '  -1.11%  '
$ws.Range('E30').Value = '  +2.80%  '
$ws.Range('D31').Value = "'27.56"
$ws.Range('E31').Value = '  +0.35%  '
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('E33').Value = '  -3.95%  '
$ws.Range('E34').Value = '  -4.47%  '
$ws.Range('E35').Value = '  -1.18%  '
$ws.Range('E36').Value = '  -2.13%  '
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('E38').Value = '  -4.53%  '
$ws.Range('D40').Value = "'423.02"
$ws.Range('E40').Value = '  -1.91%  '
$ws.Range('D41').Value = "'8.39"
$ws.Range('E41').Value = '  +0.13%  '
$ws.Range('D42').Value = '2.972.13'
$ws.Range('E42').Value = '  +2.04%  '
$ws.Range('D43').Value = "'2.75"
$ws.Range('E43').Value = '  -5.92%  '
$ws.Range('E44').Value = '  -8.13%  '
$ws.Range('D45').Value = "'0.266"
$ws.Range('E45').Value = '  +2.44%  '
$ws.Range('E46').Value = '  -0.89%  '
$ws.Range('E48').Value = '  +0.49%  '
$ws.Range('E49').Value = '  -3.16%  '
$ws.Range('E50').Value = '  -0.53%  '
$ws.Range('D51').Value = "'121.49"
$ws.Range('E51').Value = '  +0.41%  '
